$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" and "Volume(1h)" columns store numeric-looking values as plain
# text (inline strings) in the workbook. Setting NumberFormat to Text ("@")
# before assigning the new value keeps Excel from auto-converting strings
# like "305.34" or "2.26%" into real numbers/percentages.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "305.34"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "2.26%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "31.73"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "-0.17%"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.174"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "1.43%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07529"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "2.331"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "32.19%"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "8.018"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "3.46%"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9156"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "-1.36%"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1742"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "1.93%"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07571"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "3.94%"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.08262"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "3.96%"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03032"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-0.84%"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09942"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "0.42%"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001504"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "0.95%"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.006109"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "-6.53%"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.502"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1.50%"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.880"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "2.29%"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "0.67%"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-0.49%"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.1337"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "0.94%"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.652"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "2.03%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.04619"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "-0.69%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.1562"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "0.88%"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001261"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "3.64%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004540"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "2.63%"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-7.26%"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0002735"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "47.24%"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01783"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "5.83%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04596"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "0.88%"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.007303"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "3.49%"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1366"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "3.01%"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002195"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "6.66%"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-15.76%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00006502"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "7.01%"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-57.48%"
